$wb = $excel.ActiveWorkbook

# ========================================================
# Section_A sheet updates (timetable view with room numbers)
# ========================================================
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B2").Value = "MA263 [C302]"
$ws.Range("C2").Value = "EC261 [C404]"
$ws.Range("D2").Value = "EC261 [C404]"
$ws.Range("E2").Value = "EC263 [C002]"
$ws.Range("F2").Value = "CS307 [C302]"
$ws.Range("B3").Value = "MA262 [C003]"
$ws.Range("C3").Value = "EC262 [C204]"
$ws.Range("D3").Value = "EC262 [C204]"
$ws.Range("E3").Value = "MA261 [C205]"
$ws.Range("F3").Value = "EC263 [C002]"
$ws.Range("B5").Value = "ELECTIVE_B3 [C004]"
$ws.Range("C5").Value = "CS307 [C302]"
$ws.Range("D5").Value = "ELECTIVE_B3 [C004]"
$ws.Range("E5").Value = "EC262 (Lab) [L408]"
$ws.Range("F5").Value = "EC263 (Lab) [L408]"
$ws.Range("C6").Value = "MA261 (Tutorial) [C002]"
$ws.Range("D6").Value = "CS307 (Tutorial) [C203]"
$ws.Range("E6").Value = "EC262 (Lab) [L408]"
$ws.Range("F6").Value = "EC263 (Lab) [L408]"
$ws.Range("B7").Value = "MA261 [C205]"
$ws.Range("C7").Value = "MA263 [C302]"
$ws.Range("D7").Value = "Free"
$ws.Range("E7").Value = "MA262 [C003]"
$ws.Range("F7").Value = "Free"
$ws.Range("B8").Value = "Free"
$ws.Range("D8").Value = "MA262 (Tutorial) [C202]"
$ws.Range("E8").Value = "EC261 (Tutorial) [C305]"

# ========================================================
# Section_B sheet updates (timetable view with room numbers)
# ========================================================
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("C2").Value = "MA263 [C304]"
$ws.Range("D2").Value = "EC261 [C104]"
$ws.Range("E2").Value = "EC263 [C203]"
$ws.Range("F2").Value = "EC263 [C203]"
$ws.Range("B3").Value = "EC262 [C203]"
$ws.Range("C3").Value = "Free"
$ws.Range("D3").Value = "CS307 [C201]"
$ws.Range("E3").Value = "EC262 [C203]"
$ws.Range("F3").Value = "MA262 [C101]"
$ws.Range("B5").Value = "ELECTIVE_B3 [C304]"
$ws.Range("C5").Value = "EC261 [C104]"
$ws.Range("D5").Value = "ELECTIVE_B3 [C304]"
$ws.Range("E5").Value = "MA262 [C101]"
$ws.Range("F5").Value = "MA261 [C203]"
$ws.Range("B6").Value = "Free"
$ws.Range("C6").Value = "MA262 (Tutorial) [C305]"
$ws.Range("D6").Value = "Free"
$ws.Range("E6").Value = "CS307 (Tutorial) [C003]"
$ws.Range("F6").Value = "Free"
$ws.Range("B7").Value = "MA263 [C304]"
$ws.Range("C7").Value = "EC263 (Lab) [L406]"
$ws.Range("D7").Value = "EC262 (Lab) [L408]"
$ws.Range("E7").Value = "MA261 [C203]"
$ws.Range("F7").Value = "CS307 [C201]"
$ws.Range("B8").Value = "MA261 (Tutorial) [C302]"
$ws.Range("C8").Value = "EC263 (Lab) [L406]"
$ws.Range("D8").Value = "EC262 (Lab) [L408]"
$ws.Range("E8").Value = "Free"
$ws.Range("F8").Value = "EC261 (Tutorial) [C205]"

# ========================================================
# Classroom_Utilization sheet updates (numeric hours/rates)
# ========================================================
$ws = $wb.Worksheets.Item("Classroom_Utilization")
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 0.8
$ws.Range("G3").Value = 10
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 0.8
$ws.Range("G4").Value = 10
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0.6
$ws.Range("G6").Value = 7.5
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 0.6
$ws.Range("G9").Value = 7.5
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 0.6
$ws.Range("G13").Value = 7.5
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.2
$ws.Range("G14").Value = 2.5
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 25
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 1.4
$ws.Range("G22").Value = 17.5
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 1.2
$ws.Range("G24").Value = 15
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0.4
$ws.Range("G25").Value = 5
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = 0.6
$ws.Range("G32").Value = 7.5
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("D34").Value = 2.5
$ws.Range("E34").Value = 0.5
$ws.Range("G34").Value = 6.25
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("D36").Value = 7.5
$ws.Range("E36").Value = 1.5
$ws.Range("G36").Value = 18.75

# ========================================================
# Classroom_Allocation sheet updates (flat allocation table)
# ========================================================
$ws = $wb.Worksheets.Item("Classroom_Allocation")
$ws.Range("F2").Value = "MA263"
$ws.Range("G2").Value = "C302"
$ws.Range("H2").Value = "classroom"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "96"
$ws.Range("F3").Value = "MA262"
$ws.Range("G3").Value = "C003"
$ws.Range("H3").Value = "large classroom"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "135"
$ws.Range("G4").Value = "C004"
$ws.Range("H4").Value = "Auditorium"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "240"
$ws.Range("F5").Value = "MA261"
$ws.Range("D6").Value = "Tue"
$ws.Range("E6").Value = "09:00-10:30"
$ws.Range("F6").Value = "EC261"
$ws.Range("G6").Value = "C404"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "78"
$ws.Range("E7").Value = "10:30-12:00"
$ws.Range("F7").Value = "EC262"
$ws.Range("G7").Value = "C204"
$ws.Range("E8").Value = "13:00-14:30"
$ws.Range("F8").Value = "CS307"
$ws.Range("G8").Value = "C302"
$ws.Range("H8").Value = "classroom"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "96"
$ws.Range("E9").Value = "14:30-15:30"
$ws.Range("F9").Value = "MA261 (Tutorial)"
$ws.Range("G9").Value = "C002"
$ws.Range("H9").Value = "large classroom"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "116"
$ws.Range("E10").Value = "15:30-17:00"
$ws.Range("F10").Value = "MA263"
$ws.Range("G10").Value = "C302"
$ws.Range("D11").Value = "Wed"
$ws.Range("E11").Value = "09:00-10:30"
$ws.Range("F11").Value = "EC261"
$ws.Range("G11").Value = "C404"
$ws.Range("H11").Value = "classroom"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "78"
$ws.Range("E12").Value = "10:30-12:00"
$ws.Range("F12").Value = "EC262"
$ws.Range("G12").Value = "C204"
$ws.Range("G13").Value = "C004"
$ws.Range("H13").Value = "Auditorium"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "240"
$ws.Range("F14").Value = "CS307 (Tutorial)"
$ws.Range("G14").Value = "C203"
$ws.Range("E15").Value = "17:00-18:00"
$ws.Range("F15").Value = "MA262 (Tutorial)"
$ws.Range("G15").Value = "C202"
$ws.Range("D16").Value = "Thu"
$ws.Range("E16").Value = "09:00-10:30"
$ws.Range("F16").Value = "EC263"
$ws.Range("G16").Value = "C002"
$ws.Range("H16").Value = "large classroom"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "116"
$ws.Range("F17").Value = "MA261"
$ws.Range("G17").Value = "C205"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "96"
$ws.Range("G18").Value = "L408"
$ws.Range("H18").Value = "classroom without projector"
$ws.Range("G19").Value = "L408"
$ws.Range("H19").Value = "classroom without projector"
$ws.Range("F20").Value = "MA262"
$ws.Range("D21").Value = "Thu"
$ws.Range("E21").Value = "17:00-18:00"
$ws.Range("F21").Value = "EC261 (Tutorial)"
$ws.Range("G21").Value = "C305"
$ws.Range("H21").Value = "classroom"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "96"
$ws.Range("E22").Value = "09:00-10:30"
$ws.Range("F22").Value = "CS307"
$ws.Range("G22").Value = "C302"
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "96"
$ws.Range("E23").Value = "10:30-12:00"
$ws.Range("F23").Value = "EC263"
$ws.Range("G23").Value = "C002"
$ws.Range("H23").Value = "large classroom"
$ws.Range("I23").NumberFormat = "@"
$ws.Range("I23").Value = "116"
$ws.Range("J23").Value = "Projector"
$ws.Range("E24").Value = "13:00-14:30"
$ws.Range("G24").Value = "L408"
$ws.Range("H24").Value = "classroom without projector"
$ws.Range("E25").Value = "14:30-15:30"
$ws.Range("F25").Value = "EC263 (Lab)"
$ws.Range("G25").Value = "L408"
$ws.Range("H25").Value = "classroom without projector"
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value = "78"
$ws.Range("J25").Value = "Computers"
$ws.Range("G26").Value = "C203"
$ws.Range("G27").Value = "C304"
$ws.Range("E28").Value = "15:30-17:00"
$ws.Range("F28").Value = "MA263"
$ws.Range("G28").Value = "C304"
$ws.Range("I28").NumberFormat = "@"
$ws.Range("I28").Value = "96"
$ws.Range("E29").Value = "17:00-18:00"
$ws.Range("F29").Value = "MA261 (Tutorial)"
$ws.Range("G29").Value = "C302"
$ws.Range("D30").Value = "Tue"
$ws.Range("E30").Value = "09:00-10:30"
$ws.Range("F30").Value = "MA263"
$ws.Range("G30").Value = "C304"
$ws.Range("E31").Value = "13:00-14:30"
$ws.Range("F31").Value = "EC261"
$ws.Range("G31").Value = "C104"
$ws.Range("E32").Value = "14:30-15:30"
$ws.Range("F32").Value = "MA262 (Tutorial)"
$ws.Range("G32").Value = "C305"
$ws.Range("E33").Value = "15:30-17:00"
$ws.Range("F33").Value = "EC263 (Lab)"
$ws.Range("G33").Value = "L406"
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = "78"
$ws.Range("J33").Value = "Computers"
$ws.Range("E34").Value = "17:00-18:00"
$ws.Range("F34").Value = "EC263 (Lab)"
$ws.Range("G34").Value = "L406"
$ws.Range("D35").Value = "Wed"
$ws.Range("E35").Value = "09:00-10:30"
$ws.Range("F35").Value = "EC261"
$ws.Range("G35").Value = "C104"
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = "96"
$ws.Range("J35").Value = "Projector"
$ws.Range("G36").Value = "C201"
$ws.Range("G37").Value = "C304"
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "96"
$ws.Range("E38").Value = "15:30-17:00"
$ws.Range("F38").Value = "EC262 (Lab)"
$ws.Range("G38").Value = "L408"
$ws.Range("H38").Value = "classroom without projector"
$ws.Range("J38").Value = "Computers"
$ws.Range("E39").Value = "17:00-18:00"
$ws.Range("F39").Value = "EC262 (Lab)"
$ws.Range("G39").Value = "L408"
$ws.Range("H39").Value = "classroom without projector"
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "78"
$ws.Range("J39").Value = "Computers"
$ws.Range("F40").Value = "EC263"
$ws.Range("F41").Value = "EC262"
$ws.Range("G41").Value = "C203"
$ws.Range("F42").Value = "MA262"
$ws.Range("G42").Value = "C101"
$ws.Range("E43").Value = "14:30-15:30"
$ws.Range("F43").Value = "CS307 (Tutorial)"
$ws.Range("G43").Value = "C003"
$ws.Range("H43").Value = "large classroom"
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = "135"
$ws.Range("D44").Value = "Thu"
$ws.Range("E44").Value = "15:30-17:00"
$ws.Range("G44").Value = "C203"
$ws.Range("E45").Value = "09:00-10:30"
$ws.Range("F45").Value = "EC263"
$ws.Range("G45").Value = "C203"
$ws.Range("E46").Value = "10:30-12:00"
$ws.Range("F46").Value = "MA262"
$ws.Range("G46").Value = "C101"
$ws.Range("J46").Value = "Projector"
$ws.Range("E47").Value = "13:00-14:30"
$ws.Range("F47").Value = "MA261"
$ws.Range("G47").Value = "C203"
$ws.Range("J47").Value = "Projector"
$ws.Range("F48").Value = "CS307"
$ws.Range("G48").Value = "C201"
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value = "96"
$ws.Range("F49").Value = "EC261 (Tutorial)"
$ws.Range("G49").Value = "C205"
